$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected x_m (E) / y_m (F) sensor coordinates ("Ya funcionan las distancias, queda mejorar el modelo").
$ws.Cells.Item(2, 5).Value = 26.46933819339602
$ws.Cells.Item(2, 6).Value = 2.809523913471157
$ws.Cells.Item(3, 5).Value = 26.46933819339602
$ws.Cells.Item(3, 6).Value = 2.809523913471157
$ws.Cells.Item(4, 5).Value = 26.46933819339602
$ws.Cells.Item(4, 6).Value = 2.809523913471157
$ws.Cells.Item(5, 5).Value = 26.46933819339602
$ws.Cells.Item(5, 6).Value = 2.809523913471157
$ws.Cells.Item(6, 5).Value = 16.207180173491
$ws.Cells.Item(6, 6).Value = 2.774837043003561
$ws.Cells.Item(7, 5).Value = 16.207180173491
$ws.Cells.Item(7, 6).Value = 2.774837043003561
$ws.Cells.Item(8, 5).Value = 18.27817929037069
$ws.Cells.Item(8, 6).Value = 4.153471912720414
$ws.Cells.Item(9, 5).Value = 18.27817929037069
$ws.Cells.Item(9, 6).Value = 4.153471912720414
$ws.Cells.Item(10, 5).Value = 18.27817929037069
$ws.Cells.Item(10, 6).Value = 4.153471912720414
$ws.Cells.Item(11, 5).Value = 18.27817929037069
$ws.Cells.Item(11, 6).Value = 4.153471912720414
$ws.Cells.Item(12, 5).Value = 20.46173983174909
$ws.Cells.Item(12, 6).Value = 1.626014051025849
$ws.Cells.Item(13, 5).Value = 20.46173983174909
$ws.Cells.Item(13, 6).Value = 1.626014051025849
$ws.Cells.Item(14, 5).Value = 20.46173983174909
$ws.Cells.Item(14, 6).Value = 1.626014051025849
$ws.Cells.Item(15, 5).Value = 20.46173983174909
$ws.Cells.Item(15, 6).Value = 1.626014051025849
$ws.Cells.Item(16, 5).Value = 19.36995956105989
$ws.Cells.Item(16, 6).Value = 2.889742981873131
$ws.Cells.Item(17, 5).Value = 19.36995956105989
$ws.Cells.Item(17, 6).Value = 2.889742981873131
$ws.Cells.Item(18, 5).Value = 19.36995956105989
$ws.Cells.Item(18, 6).Value = 2.889742981873131
$ws.Cells.Item(19, 5).Value = 19.36995956105989
$ws.Cells.Item(19, 6).Value = 2.889742981873131
$ws.Cells.Item(20, 5).Value = 19.36995956105989
$ws.Cells.Item(20, 6).Value = 2.889742981873131
$ws.Cells.Item(21, 5).Value = 19.36995956105989
$ws.Cells.Item(21, 6).Value = 2.889742981873131
$ws.Cells.Item(22, 5).Value = 22.64530037312749
$ws.Cells.Item(22, 6).Value = 4.153471912720414
$ws.Cells.Item(23, 5).Value = 22.64530037312749
$ws.Cells.Item(23, 6).Value = 4.153471912720414
$ws.Cells.Item(24, 5).Value = 22.64530037312749
$ws.Cells.Item(24, 6).Value = 4.153471912720414
$ws.Cells.Item(25, 5).Value = 22.64530037312749
$ws.Cells.Item(25, 6).Value = 4.153471912720414
$ws.Cells.Item(26, 5).Value = 24.82886091450589
$ws.Cells.Item(26, 6).Value = 1.626014051025849
$ws.Cells.Item(27, 5).Value = 24.82886091450589
$ws.Cells.Item(27, 6).Value = 1.626014051025849
$ws.Cells.Item(28, 5).Value = 24.82886091450589
$ws.Cells.Item(28, 6).Value = 1.626014051025849
$ws.Cells.Item(29, 5).Value = 24.82886091450589
$ws.Cells.Item(29, 6).Value = 1.626014051025849
$ws.Cells.Item(30, 5).Value = 23.73708064381669
$ws.Cells.Item(30, 6).Value = 2.889742981873131
$ws.Cells.Item(31, 5).Value = 23.73708064381669
$ws.Cells.Item(31, 6).Value = 2.889742981873131
$ws.Cells.Item(32, 5).Value = 23.73708064381669
$ws.Cells.Item(32, 6).Value = 2.889742981873131
$ws.Cells.Item(33, 5).Value = 23.73708064381669
$ws.Cells.Item(33, 6).Value = 2.889742981873131
$ws.Cells.Item(34, 5).Value = 23.73708064381669
$ws.Cells.Item(34, 6).Value = 2.889742981873131
$ws.Cells.Item(35, 5).Value = 23.73708064381669
$ws.Cells.Item(35, 6).Value = 2.889742981873131
$ws.Cells.Item(36, 5).Value = 0.3923531079605417
$ws.Cells.Item(36, 6).Value = 2.816622836057509
$ws.Cells.Item(37, 5).Value = 0.3923531079605417
$ws.Cells.Item(37, 6).Value = 2.816622836057509
$ws.Cells.Item(38, 5).Value = 0.3923531079605417
$ws.Cells.Item(38, 6).Value = 2.816622836057509
$ws.Cells.Item(39, 5).Value = 0.3923531079605417
$ws.Cells.Item(39, 6).Value = 2.816622836057509
$ws.Cells.Item(40, 5).Value = 6.003652662094717
$ws.Cells.Item(40, 6).Value = 4.055201217336458
$ws.Cells.Item(41, 5).Value = 6.003652662094717
$ws.Cells.Item(41, 6).Value = 4.055201217336458
$ws.Cells.Item(42, 5).Value = 6.003652662094717
$ws.Cells.Item(42, 6).Value = 4.055201217336458
$ws.Cells.Item(43, 5).Value = 6.003652662094717
$ws.Cells.Item(43, 6).Value = 4.055201217336458
$ws.Cells.Item(44, 5).Value = 7.681135109620515
$ws.Cells.Item(44, 6).Value = 1.71574406149687
$ws.Cells.Item(45, 5).Value = 7.681135109620515
$ws.Cells.Item(45, 6).Value = 1.71574406149687
$ws.Cells.Item(46, 5).Value = 7.681135109620515
$ws.Cells.Item(46, 6).Value = 1.71574406149687
$ws.Cells.Item(47, 5).Value = 7.681135109620515
$ws.Cells.Item(47, 6).Value = 1.71574406149687
$ws.Cells.Item(48, 5).Value = 6.842393885857616
$ws.Cells.Item(48, 6).Value = 2.885472639416664
$ws.Cells.Item(49, 5).Value = 6.842393885857616
$ws.Cells.Item(49, 6).Value = 2.885472639416664
$ws.Cells.Item(50, 5).Value = 6.842393885857616
$ws.Cells.Item(50, 6).Value = 2.885472639416664
$ws.Cells.Item(51, 5).Value = 6.842393885857616
$ws.Cells.Item(51, 6).Value = 2.885472639416664
$ws.Cells.Item(52, 5).Value = 6.842393885857616
$ws.Cells.Item(52, 6).Value = 2.885472639416664
$ws.Cells.Item(53, 5).Value = 6.842393885857616
$ws.Cells.Item(53, 6).Value = 2.885472639416664
$ws.Cells.Item(54, 5).Value = 9.358617557146312
$ws.Cells.Item(54, 6).Value = 4.055201217336458
$ws.Cells.Item(55, 5).Value = 9.358617557146312
$ws.Cells.Item(55, 6).Value = 4.055201217336458
$ws.Cells.Item(56, 5).Value = 9.358617557146312
$ws.Cells.Item(56, 6).Value = 4.055201217336458
$ws.Cells.Item(57, 5).Value = 9.358617557146312
$ws.Cells.Item(57, 6).Value = 4.055201217336458
$ws.Cells.Item(58, 5).Value = 11.03610000467211
$ws.Cells.Item(58, 6).Value = 1.71574406149687
$ws.Cells.Item(59, 5).Value = 11.03610000467211
$ws.Cells.Item(59, 6).Value = 1.71574406149687
$ws.Cells.Item(60, 5).Value = 11.03610000467211
$ws.Cells.Item(60, 6).Value = 1.71574406149687
$ws.Cells.Item(61, 5).Value = 11.03610000467211
$ws.Cells.Item(61, 6).Value = 1.71574406149687
$ws.Cells.Item(62, 5).Value = 10.19735878090921
$ws.Cells.Item(62, 6).Value = 2.885472639416664
$ws.Cells.Item(63, 5).Value = 10.19735878090921
$ws.Cells.Item(63, 6).Value = 2.885472639416664
$ws.Cells.Item(64, 5).Value = 10.19735878090921
$ws.Cells.Item(64, 6).Value = 2.885472639416664
$ws.Cells.Item(65, 5).Value = 10.19735878090921
$ws.Cells.Item(65, 6).Value = 2.885472639416664
$ws.Cells.Item(66, 5).Value = 10.19735878090921
$ws.Cells.Item(66, 6).Value = 2.885472639416664
$ws.Cells.Item(67, 5).Value = 10.19735878090921
$ws.Cells.Item(67, 6).Value = 2.885472639416664
$ws.Cells.Item(68, 5).Value = 12.7135824521979
$ws.Cells.Item(68, 6).Value = 4.055201217336458
$ws.Cells.Item(69, 5).Value = 12.7135824521979
$ws.Cells.Item(69, 6).Value = 4.055201217336458
$ws.Cells.Item(70, 5).Value = 12.7135824521979
$ws.Cells.Item(70, 6).Value = 4.055201217336458
$ws.Cells.Item(71, 5).Value = 12.7135824521979
$ws.Cells.Item(71, 6).Value = 4.055201217336458
$ws.Cells.Item(72, 5).Value = 14.3910648997237
$ws.Cells.Item(72, 6).Value = 1.71574406149687
$ws.Cells.Item(73, 5).Value = 14.3910648997237
$ws.Cells.Item(73, 6).Value = 1.71574406149687
$ws.Cells.Item(74, 5).Value = 14.3910648997237
$ws.Cells.Item(74, 6).Value = 1.71574406149687
$ws.Cells.Item(75, 5).Value = 14.3910648997237
$ws.Cells.Item(75, 6).Value = 1.71574406149687
$ws.Cells.Item(76, 5).Value = 13.55232367596081
$ws.Cells.Item(76, 6).Value = 2.885472639416664
$ws.Cells.Item(77, 5).Value = 13.55232367596081
$ws.Cells.Item(77, 6).Value = 2.885472639416664
$ws.Cells.Item(78, 5).Value = 13.55232367596081
$ws.Cells.Item(78, 6).Value = 2.885472639416664
$ws.Cells.Item(79, 5).Value = 13.55232367596081
$ws.Cells.Item(79, 6).Value = 2.885472639416664
$ws.Cells.Item(80, 5).Value = 13.55232367596081
$ws.Cells.Item(80, 6).Value = 2.885472639416664
$ws.Cells.Item(81, 5).Value = 13.55232367596081
$ws.Cells.Item(81, 6).Value = 2.885472639416664
$ws.Cells.Item(82, 5).Value = 2.032746778420167
$ws.Cells.Item(82, 6).Value = 1.529443221960014
$ws.Cells.Item(83, 5).Value = 2.032746778420167
$ws.Cells.Item(83, 6).Value = 1.529443221960014
$ws.Cells.Item(84, 5).Value = 15.24718637090686
$ws.Cells.Item(84, 6).Value = 2.607608394372307
$ws.Cells.Item(85, 5).Value = 15.24718637090686
$ws.Cells.Item(85, 6).Value = 2.607608394372307
$ws.Cells.Item(86, 5).Value = 26.60667347863288
$ws.Cells.Item(86, 6).Value = 2.789103681553595
$ws.Cells.Item(87, 5).Value = 26.60667347863288
$ws.Cells.Item(87, 6).Value = 2.789103681553595
$ws.Cells.Item(88, 5).Value = 26.60667347863288
$ws.Cells.Item(88, 6).Value = 2.789103681553595
$ws.Cells.Item(89, 5).Value = 26.60667347863288
$ws.Cells.Item(89, 6).Value = 2.789103681553595
$ws.Cells.Item(90, 5).Value = 11.58445098476159
$ws.Cells.Item(90, 6).Value = 2.620270498250707
$ws.Cells.Item(91, 5).Value = 11.58445098476159
$ws.Cells.Item(91, 6).Value = 2.620270498250707
